# MapaConceptual.pptx - "mapa conceptual matemáticas 9 tema 2"
#
# Applies the authored changes to the single concept-map slide:
#   - nudge the "La adición" node slightly to the left
#   - collapse several runs that were split mid-word/mid-phrase back into
#     single runs (no visible text change, just how PowerPoint re-flowed
#     the runs when the text was retouched)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Id -eq $id) {
            return $shape
        }
    }
    return $null
}

# Merge a run of characters back into a single run by re-writing exactly
# the same visible text through a real (non no-op) TextRange mutation, so
# PowerPoint-style run-coalescing kicks in and keeps the first run's
# formatting.
function Set-MergedText($textRange, $start, $length, $text) {
    $sub = $textRange.Characters($start, $length)
    $sub.Text = $text
}

# --- "Rectángulo 4" / Nodo01 ("La adición") : reposition ------------------
$nodo01 = Get-ShapeById $s.Shapes 5
$nodo01.Left = 80.0816
$nodo01.Top = 83.0137

# --- "CuadroTexto 129" ("las básicas son") --------------------------------
$shp130 = Get-ShapeById $s.Shapes 130
$tr130 = $shp130.TextFrame.TextRange
Set-MergedText $tr130 1 $tr130.Length "las básicas son"

# --- "Rectángulo 168" (bullet list ending in "opuesto aditivo") ----------
$shp169 = Get-ShapeById $s.Shapes 169
$tr169 = $shp169.TextFrame.TextRange
$para5 = $tr169.Paragraphs(5, 1)
Set-MergedText $tr169 $para5.Start $para5.Length "opuesto aditivo"

# --- "Rectángulo 122" ("se adicionan ...") --------------------------------
$shp123 = Get-ShapeById $s.Shapes 123
$tr123 = $shp123.TextFrame.TextRange
Set-MergedText $tr123 1 13 "se adicionan "

# --- "Rectángulo 176" ("al sustraer dos enteros se obtiene un entero") ---
$shp177 = Get-ShapeById $s.Shapes 177
$tr177 = $shp177.TextFrame.TextRange
Set-MergedText $tr177 1 $tr177.Length "al sustraer dos enteros se obtiene un entero"

# --- "CuadroTexto 182" ("satisface únicamente la") ------------------------
$shp183 = Get-ShapeById $s.Shapes 183
$tr183 = $shp183.TextFrame.TextRange
Set-MergedText $tr183 1 $tr183.Length "satisface únicamente la"

# --- "CuadroTexto 240" ("se escribe") -------------------------------------
$shp241 = Get-ShapeById $s.Shapes 241
$tr241 = $shp241.TextFrame.TextRange
Set-MergedText $tr241 1 $tr241.Length "se escribe"

# --- "Rectángulo 243" ("cociente exacto de los valores absolutos de los
#      números") : merge the first three runs, and the last two runs,
#      leaving "los valores " (the 4th run) untouched ---------------------
$shp244 = Get-ShapeById $s.Shapes 244
$tr244 = $shp244.TextFrame.TextRange
Set-MergedText $tr244 1 19 "cociente exacto de "
Set-MergedText $tr244 32 24 "absolutos de los números"

# --- "CuadroTexto 93" ("corresponde al") ----------------------------------
$shp94 = Get-ShapeById $s.Shapes 94
$tr94 = $shp94.TextFrame.TextRange
Set-MergedText $tr94 1 $tr94.Length "corresponde al"
